$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verwendete Lizenzen")

# Remove the "ZXing" row (row 6) entirely - shifts subsequent rows up
$ws.Rows.Item(6).Delete()
